$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = 0.95432611918374544
$ws.Range("T1").Value = 0.88460935080061431
$ws.Range("AM1").Value = 0.94964801443745217
$ws.Range("AP2").Value = 0.72434425643317613
$ws.Range("I3").Value = 0.76573902215764755
$ws.Range("AE3").Value = 0.85219670907259149
$ws.Range("F4").Value = 0.76788849283222371
$ws.Range("AS4").Value = 0.76356437512763597
$ws.Range("C5").Value = 0.98872631165080649
$ws.Range("AW6").Value = 0.91576598540248444
$ws.Range("I7").Value = 0.97211331198335449
$ws.Range("Z8").Value = 0.80468935033190736
$ws.Range("E9").Value = 0.82631453270331945
$ws.Range("AN9").Value = 0.99001735557272119
$ws.Range("C10").Value = 0.84639842120813868
$ws.Range("M11").Value = 0.84009658287327449
$ws.Range("BO11").Value = 0.75184774385446007
$ws.Range("AD12").Value = 0.64752655964938821
$ws.Range("AW12").Value = 0.78082926713317202
$ws.Range("AN13").Value = 0.99196762976849018
$ws.Range("BL13").Value = 0.97223412040364077
$ws.Range("AD14").Value = 0.69921938629675529
$ws.Range("N16").Value = 0.55162320250706542
$ws.Range("O16").Value = 0.92644723967936016
$ws.Range("P17").Value = 0.68337862758266299
$ws.Range("T17").Value = 0.72809862439394224
$ws.Range("BK18").Value = 0.91158302311635764
$ws.Range("F19").Value = 0.77123631124136538
$ws.Range("AU19").Value = 0.90853234309724384
$ws.Range("AX19").Value = 0.85139183951226527
$ws.Range("G21").Value = 0.6932506581952893
$ws.Range("AW21").Value = 0.86737929890488541
$ws.Range("F22").Value = 0.87283035808039089
$ws.Range("P22").Value = 0.73774743812973864
$ws.Range("W22").Value = 0.84606833220088307
$ws.Range("J23").Value = 0.62304717349690752
$ws.Range("X23").Value = 0.65470675049781757
$ws.Range("J24").Value = 0.8812108029229162
$ws.Range("AI24").Value = 0.99776715764680035
$ws.Range("P25").Value = 0.80895931210227601
$ws.Range("BG25").Value = 0.9618151886689541
$ws.Range("O26").Value = 0.92085189824502345
$ws.Range("AX26").Value = 0.58700164610446448
$ws.Range("AZ26").Value = 0.75486424644941885
$ws.Range("AC27").Value = 0.86096872902376598
$ws.Range("AZ27").Value = 0.83128106587302775
$ws.Range("BK27").Value = 0.84489171872988855
$ws.Range("AC28").Value = 0.72005052486207799
$ws.Range("AJ28").Value = 0.69320726816192657
$ws.Range("I29").Value = 0.97867620841380454
$ws.Range("AR29").Value = 0.99602719477268087
$ws.Range("E30").Value = 0.9532358525631992
$ws.Range("AC30").Value = 0.94814254677650334
$ws.Range("AE30").Value = 0.63023648490860895
$ws.Range("H31").Value = 0.90727483359059646
$ws.Range("Q32").Value = 0.95751675848558926
$ws.Range("W32").Value = 0.99465393890325116
$ws.Range("AN32").Value = 0.9085236072349443
$ws.Range("Y33").Value = 0.96421313908658046
$ws.Range("AB33").Value = 0.7970259513049508
$ws.Range("AS33").Value = 0.8339017786411298
$ws.Range("AW33").Value = 0.81126595415019054
$ws.Range("AX33").Value = 0.88743527958539836
$ws.Range("Q34").Value = 0.98446862006646008
$ws.Range("BD34").Value = 0.78084961434370082
$ws.Range("AD35").Value = 0.75029172521542542
$ws.Range("AH35").Value = 0.59203241093550041
$ws.Range("AJ35").Value = 0.77799427694979439
$ws.Range("K36").Value = 0.96293807520532981
$ws.Range("BA36").Value = 0.91209660073059218
$ws.Range("BH36").Value = 0.78686695959748143
$ws.Range("G37").Value = 0.90915414831715813
$ws.Range("J37").Value = 0.80200350867367787
$ws.Range("P37").Value = 0.87457216625218548
$ws.Range("AC38").Value = 0.92215552635683595
$ws.Range("AB39").Value = 0.83218302212732731
$ws.Range("BH39").Value = 0.69321619147968505
$ws.Range("O40").Value = 0.94744645968207619
$ws.Range("AL40").Value = 0.85396047307026657
$ws.Range("AO40").Value = 0.67496930016656065
$ws.Range("AP41").Value = 0.70907233148909021
$ws.Range("X42").Value = 0.92423753469222336
$ws.Range("BB42").Value = 0.93032152156488879
$ws.Range("AR43").Value = 0.86748770854836466
$ws.Range("AS43").Value = 0.94347068263374778
$ws.Range("AH44").Value = 0.65204169696461411
$ws.Range("AZ44").Value = 0.68656398843399968
$ws.Range("AC46").Value = 0.86261448725391032
$ws.Range("AV46").Value = 0.82776495157700158
$ws.Range("O47").Value = 0.96428849580111942
$ws.Range("AT47").Value = 0.86092624367333825
$ws.Range("AV47").Value = 0.76928140146789015
$ws.Range("V48").Value = 0.99875248301392427
$ws.Range("AK48").Value = 0.9450460472982698
$ws.Range("K49").Value = 0.71410050684448478
$ws.Range("BG50").Value = 0.91010585439323122
$ws.Range("BA51").Value = 0.86007076039721531
$ws.Range("BP51").Value = 0.77724199008259676
$ws.Range("AP52").Value = 0.82296928651590784
$ws.Range("AY52").Value = 0.97161593620943854
$ws.Range("D53").Value = 0.80960473780151465
$ws.Range("AV53").Value = 0.82113733629612584
$ws.Range("AI54").Value = 0.72320167699428706
$ws.Range("AZ54").Value = 0.68077342867986146
$ws.Range("BD54").Value = 0.96187464986076876
$ws.Range("R55").Value = 0.8579611018568154
$ws.Range("T55").Value = 0.98456582364649303
$ws.Range("AU55").Value = 0.69997825612645426
$ws.Range("BL55").Value = 0.93458143687797635
$ws.Range("M56").Value = 0.86647344159727013
$ws.Range("S56").Value = 0.79496210096548814
$ws.Range("BI56").Value = 0.63914242463069959
$ws.Range("L57").Value = 0.77014925612866947
$ws.Range("AT57").Value = 0.69351537822179077
$ws.Range("BC57").Value = 0.66365603848899213
$ws.Range("BD57").Value = 0.89703427041748141
$ws.Range("BI57").Value = 0.89178085871367208
$ws.Range("B59").Value = 0.7587959292680293
$ws.Range("BF59").Value = 0.97237563584530284
$ws.Range("BF60").Value = 0.83341797965962661
$ws.Range("X62").Value = 0.85885745508704625
$ws.Range("AY62").Value = 0.73949636155201781
$ws.Range("BH62").Value = 0.7822506207830362
$ws.Range("H63").Value = 0.72106588715193887
$ws.Range("BB64").Value = 0.94697953114888034
$ws.Range("Z65").Value = 0.6870347026690109
$ws.Range("AA65").Value = 0.89614118897435957
$ws.Range("BK65").Value = 0.86671066479311587
$ws.Range("L66").Value = 0.89793414724999376
$ws.Range("BD66").Value = 0.81287637938048407
$ws.Range("BG66").Value = 0.91250989924892112
$ws.Range("BM66").Value = 0.98027831995946357
$ws.Range("W68").Value = 0.65142900559358496
$ws.Range("Y68").Value = 0.81045115419118141
$ws.Range("BM68").Value = 0.84578625090432569
$ws.Range("BO68").Value = 0.75029005898880796
